# Game of Turns - remove the high-level algorithm slide and let students
# work out their own solution.

$p = $ppt.ActivePresentation

# --- Slide 2 ("How does the game work?"): refine the problem statement and
#     add a concrete example of the turn order, moved up from the slide
#     that is being removed below. ---
$s2 = $p.Slides.Item(2)
$body = $s2.Shapes.Item(2)
$lineBreak = [char]11
$body.TextFrame.TextRange.Text = "We would like to make the signal travel among all the processes, stating from the parent process, in a orderly circular fashion for a total of M rounds.`rExample: " + $lineBreak + "P -> C1 -> C2 -> C3 -> P"

# --- Remove the old "How does the game work?" slide that spelled out the
#     step-by-step algorithm (slide 3). The following slide ("Signal")
#     shifts up to take its place. ---
$p.Slides.Item(3).Delete()
